$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Results")

$ws.Range("C16").Value = "C:\Users\COCO\onnxruntime_training_cuda_python\orttraining\orttraining\python\orttraining_pybind_state.cc:621 onnxruntime::python::addObjectMethodsForTraining::<lambda_6dd399ad6691adab5d0e0423ed8ce22d>::operator () [ONNXRuntimeError] : 1 : FAIL : Type Error: Type parameter (T) of Optype (Sub) bound to different types (tensor(float) and tensor(double) in node (onnx::Pow::53_Grad/Sub_1).`n"
$ws.Rows.Item(16).AutoFit()

$ws.Range("C39").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::172): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C51").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::262): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C52").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::264): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C64").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::317): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C65").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::319): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C69").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::337): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C70").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::339): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C93").Value = "[ShapeInferenceError] (op_type:Sub, node name: onnx::Sub::435): A typestr: T, has unsupported type: tensor(bool)"

$ws.Range("C102").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::472): X typestr: T, has unsupported type: tensor(uint8)"

$ws.Range("C223").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::922): X typestr: T, has unsupported type: tensor(uint8)"

$ws.Range("C239").Value = "[ShapeInferenceError] (op_type:Pow, node name: onnx::Pow::927): X typestr: T, has unsupported type: tensor(uint8)"
